$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65:140 down to 66:141
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new price record
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44494
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101007
$ws.Range("J65").Value = "Kiwi"
$ws.Range("K65").Value = "Hayward"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 400
$ws.Range("N65").Value = 15000
$ws.Range("O65").Value = 16000
$ws.Range("P65").Value = 15500
$ws.Range("Q65").Value = "$/caja 15 kilos"
$ws.Range("R65").Value = "Provincia de Curicó"
$ws.Range("S65").Value = 1033
$ws.Range("T65").Value = 15
